$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "actual start/end date" for the Undo feature task (row 11)
$ws.Range("H11").Value = "22 tháng 10"
$ws.Range("I11").Value = "23 tháng 10"

# Update the active selection / view to match the saved state (scrolled right, I11 selected)
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("I11").Select()
